$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.804.07'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.585.53'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'209.17"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = "'1.00"
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -3.50%  '
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").Value = "'0.0616"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = "'18.03"
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").Value = '1.805.53'
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("D13").Value = '1.576.38'
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("E15").Value = '  -2.32%  '
$ws.Range("D16").Value = '25.793.25'
$ws.Range("D17").Value = '0.0₃0721'
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = "'59.85"
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = "'191.75"
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").Value = "'5.92"
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").Value = "'141.51"
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("E30").Value = '  -5.63%  '
$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("E33").Value = '  -2.35%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").Value = '1.097.56'
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  -1.89%  '
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = '  -7.62%  '
$ws.Range("D42").Value = "'0.817"
$ws.Range("E42").Value = '  +9.37%  '
$ws.Range("E43").Value = '  +2.59%  '
$ws.Range("D44").Value = "'93.72"
$ws.Range("E44").Value = '  -4.20%  '
$ws.Range("D45").Value = '1.719.79'
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = "'53.14"
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  +0.03%  '

"Done"